{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The document ends with three empty trailing paragraphs (right before the\n// end of the body / sectPr). This change fills the *second-to-last* of\n// those empty paragraphs with the first line of new notes text, then\n// inserts a further block of paragraphs right after it (the *last* empty\n// paragraph is left untouched at the very end of the document).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\n// Second-to-last paragraph in the document \u2014 this is the empty paragraph\n// that receives the first line of the new text.\nconst target = paragraphs.items[count - 2];\n\n// The lines of text to add, in document order. Each entry is either a\n// plain string (a paragraph with a single run) or an array of strings\n// (a paragraph whose text is split across multiple runs).\nconst lines = [\n  \"You should give some thought on names for usernames, what I mean here is you can\",\n  \"have some kind of format suchas WGD0001 which is my initials ... yours KRD0001\",\n  \"so if someone has same initials just increment the number KRD0002. Another option\",\n  [\n    \"is to incorporate their dept, i.e \",\n    \"KevinFinance, SimonService ...... or if you do not want\",\n  ],\n  \"to use peoples names use generic usernames, i.e. sales01, service02, parts01 .....\",\n  \"You can create the usernames beforehand and then assign them to people.\",\n  \"\",\n  \"You can provide me with a list of usernames and puchaser names and I can load them\",\n  \"into the database so that you don\\u2019t have to type them ....\",\n];\n\nfunction firstRunText(line) {\n  return Array.isArray(line) ? line[0] : line;\n}\n\n// Fill the existing empty paragraph with the first line's (first run's)\n// text instead of inserting a brand-new paragraph for it.\ntarget.insertText(firstRunText(lines[0]), \"Replace\");\nif (Array.isArray(lines[0])) {\n  let r = target.getRange(\"End\");\n  for (let j = 1; j < lines[0].length; j++) {\n    r = r.insertText(lines[0][j], \"End\");\n  }\n}\n\n// Insert the remaining lines as new paragraphs, each right after the\n// previous one, preserving order.\nlet cursor = target;\nfor (let i = 1; i < lines.length; i++) {\n  const line = lines[i];\n  const text = firstRunText(line);\n  cursor = cursor.insertParagraph(text, \"After\");\n  if (Array.isArray(line)) {\n    let r = cursor.getRange(\"End\");\n    for (let j = 1; j < line.length; j++) {\n      r = r.insertText(line[j], \"End\");\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# The document ends with three empty trailing paragraphs (right before the\n# end of the body). This change fills the *second-to-last* of those empty\n# paragraphs with the first line of new notes text, then inserts a further\n# block of paragraphs right after it (the *last* empty paragraph is left\n# untouched at the very end of the document).\n\n$d = $word.ActiveDocument\n\n# Lines of text to add, in document order. Each entry is either a plain\n# string (single run) or an array of strings (paragraph text split across\n# multiple runs).\n$lines = @(\n    \"You should give some thought on names for usernames, what I mean here is you can\",\n    \"have some kind of format suchas WGD0001 which is my initials ... yours KRD0001\",\n    \"so if someone has same initials just increment the number KRD0002. Another option\",\n    @(\"is to incorporate their dept, i.e \", \"KevinFinance, SimonService ...... or if you do not want\"),\n    \"to use peoples names use generic usernames, i.e. sales01, service02, parts01 .....\",\n    \"You can create the usernames beforehand and then assign them to people.\",\n    \"\",\n    \"You can provide me with a list of usernames and puchaser names and I can load them\",\n    \"into the database so that you don\u2019t have to type them ....\"\n)\n\n# Second-to-last paragraph in the document \u2014 the empty paragraph that\n# receives the first line's text.\n$count = $d.Paragraphs.Count\n$target = $d.Paragraphs($count - 1)\n\n$first = $lines[0]\nif ($first -is [array]) {\n    $target.Range.Text = $first[0]\n    for ($j = 1; $j -lt $first.Count; $j++) {\n        $target.Range.InsertAfter($first[$j])\n    }\n} elseif ($first -ne \"\") {\n    $target.Range.Text = $first\n}\n\n# Insert the remaining lines as new paragraphs, each right after the\n# previous one, preserving order. A brand-new paragraph from\n# InsertParagraphAfter() already has an empty run, so a blank line needs no\n# further write (assigning Text = \"\" would instead leave a stray empty\n# <w:t/> element behind).\n$cur = $target\nfor ($i = 1; $i -lt $lines.Count; $i++) {\n    $cur.Range.InsertParagraphAfter()\n    $cur = $cur.Next()\n    $line = $lines[$i]\n    if ($line -is [array]) {\n        $cur.Range.Text = $line[0]\n        for ($j = 1; $j -lt $line.Count; $j++) {\n            $cur.Range.InsertAfter($line[$j])\n        }\n    } elseif ($line -ne \"\") {\n        $cur.Range.Text = $line\n    }\n}\n"}
